$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = 0.2479090405365451
$ws.Range("D2").Value = 0.8065040260802943

# Row 3
$ws.Range("C3").Value = 0.2233914028418744
$ws.Range("D3").Value = 0.8252921731356559

# Row 4
$ws.Range("C4").Value = -1.10784154189995
$ws.Range("D4").Value = 0.2798957992306943

# Row 5
$ws.Range("C5").Value = -2.135513497381739
$ws.Range("D5").Value = 0.04409993457764472

# Row 6
$ws.Range("C6").Value = -0.02403972993395784
$ws.Range("D6").Value = 0.9810376374017804

# Row 7
$ws.Range("C7").Value = -1.561905577968721
$ws.Range("D7").Value = 0.1325813883696225

# Row 8
$ws.Range("C8").Value = -2.590204447755234
$ws.Range("D8").Value = 0.01670384195503472

# Row 9
$ws.Range("C9").Value = -1.626746364084547
$ws.Range("D9").Value = 0.1180304866324751

# Row 10
$ws.Range("C10").Value = -2.323330950694551
$ws.Range("D10").Value = 0.02979875938686805

# Row 11
$ws.Range("C11").Value = -1.688005315932277
$ws.Range("D11").Value = 0.1055376571075366
$ws.Range("G11").Value = "No"
